# Update transition-matrix probabilities on Sheet1 per the latest simulation run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 0.2148148148148148
$ws.Range("C2").Value = 0.5148148148148148
$ws.Range("J2").Value = 0.04074074074074074
$ws.Range("P2").Value = 0.1259259259259259
$ws.Range("S2").Value = 0.1037037037037037

# Row 3
$ws.Range("B3").Value = 0.01418439716312057
$ws.Range("C3").Value = 0.007092198581560284
$ws.Range("J3").Value = 0.02127659574468085
$ws.Range("P3").Value = 0.7730496453900709
$ws.Range("S3").Value = 0.1843971631205674

# Row 4
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.3

# Row 6
$ws.Range("B6").Value = 0.07655502392344497
$ws.Range("D6").Value = 0.004784688995215311
$ws.Range("F6").Value = 0.09569377990430622
$ws.Range("J6").Value = 0.1961722488038277
$ws.Range("O6").Value = 0.03827751196172249
$ws.Range("Q6").Value = 0.138755980861244
$ws.Range("R6").Value = 0.08133971291866028
$ws.Range("S6").Value = 0.3684210526315789

# Row 7
$ws.Range("B7").Value = 0.09580838323353294
$ws.Range("D7").Value = 0.02395209580838323
$ws.Range("F7").Value = 0.04191616766467066
$ws.Range("J7").Value = 0.1197604790419162
$ws.Range("O7").Value = 0.02395209580838323
$ws.Range("Q7").Value = 0.1676646706586826
$ws.Range("R7").Value = 0.0718562874251497
$ws.Range("S7").Value = 0.4550898203592814

# Row 8
$ws.Range("B8").Value = 0.1119221411192214
$ws.Range("D8").Value = 0.024330900243309
$ws.Range("F8").Value = 0.07542579075425791
$ws.Range("J8").Value = 0.09975669099756691
$ws.Range("O8").Value = 0.0194647201946472
$ws.Range("Q8").Value = 0.1484184914841849
$ws.Range("R8").Value = 0.0827250608272506
$ws.Range("S8").Value = 0.4379562043795621

# Row 9
$ws.Range("B9").Value = 0.1237113402061856
$ws.Range("D9").Value = 0.0154639175257732
$ws.Range("F9").Value = 0.06701030927835051
$ws.Range("J9").Value = 0.08762886597938144
$ws.Range("O9").Value = 0.0154639175257732
$ws.Range("Q9").Value = 0.1804123711340206
$ws.Range("R9").Value = 0.1082474226804124
$ws.Range("S9").Value = 0.4020618556701031

# Row 10
$ws.Range("B10").Value = 0.1025390625
$ws.Range("D10").Value = 0.0234375
$ws.Range("E10").Value = 0.0009765625
$ws.Range("F10").Value = 0.080078125
$ws.Range("J10").Value = 0.099609375
$ws.Range("O10").Value = 0.013671875
$ws.Range("Q10").Value = 0.2158203125
$ws.Range("R10").Value = 0.0859375
$ws.Range("S10").Value = 0.3779296875

# Row 11
$ws.Range("G11").Value = 0.1153846153846154
$ws.Range("J11").Value = 0.06837606837606838
$ws.Range("K11").Value = 0.1581196581196581
$ws.Range("L11").Value = 0.6538461538461539
$ws.Range("S11").Value = 0.004273504273504274

# Row 12
$ws.Range("G12").Value = 0.7468354430379747
$ws.Range("J12").Value = 0.2025316455696203
$ws.Range("L12").Value = 0.02531645569620253
$ws.Range("S12").Value = 0.02531645569620253

# Row 13
$ws.Range("G13").Value = 0.6578947368421053
$ws.Range("J13").Value = 0.3421052631578947

# Row 15
$ws.Range("F15").Value = 0.0267379679144385
$ws.Range("H15").Value = 0.1390374331550802
$ws.Range("I15").Value = 0.06951871657754011
$ws.Range("J15").Value = 0.374331550802139
$ws.Range("K15").Value = 0.053475935828877
$ws.Range("M15").Value = 0.0053475935828877
$ws.Range("N15").Value = 0.0053475935828877
$ws.Range("O15").Value = 0.08021390374331551
$ws.Range("S15").Value = 0.2459893048128342

# Row 16
$ws.Range("F16").Value = 0.01219512195121951
$ws.Range("H16").Value = 0.2317073170731707
$ws.Range("I16").Value = 0.07317073170731707
$ws.Range("J16").Value = 0.3353658536585366
$ws.Range("K16").Value = 0.1280487804878049
$ws.Range("M16").Value = 0.01219512195121951
$ws.Range("O16").Value = 0.05487804878048781
$ws.Range("S16").Value = 0.1524390243902439

# Row 17
$ws.Range("F17").Value = 0.01075268817204301
$ws.Range("H17").Value = 0.1854838709677419
$ws.Range("I17").Value = 0.09946236559139784
$ws.Range("J17").Value = 0.3897849462365591
$ws.Range("K17").Value = 0.1155913978494624
$ws.Range("M17").Value = 0.02419354838709677
$ws.Range("O17").Value = 0.05913978494623656
$ws.Range("S17").Value = 0.1155913978494624

# Row 18
$ws.Range("F18").Value = 0.02906976744186046
$ws.Range("H18").Value = 0.1686046511627907
$ws.Range("I18").Value = 0.09302325581395349
$ws.Range("J18").Value = 0.436046511627907
$ws.Range("K18").Value = 0.06395348837209303
$ws.Range("M18").Value = 0.02325581395348837
$ws.Range("O18").Value = 0.06395348837209303
$ws.Range("S18").Value = 0.1220930232558139

# Row 19
$ws.Range("F19").Value = 0.008936550491510277
$ws.Range("H19").Value = 0.226988382484361
$ws.Range("I19").Value = 0.1054512957998213
$ws.Range("J19").Value = 0.3512064343163539
$ws.Range("K19").Value = 0.1000893655049151
$ws.Range("M19").Value = 0.02234137622877569
$ws.Range("N19").Value = 0.0008936550491510277
$ws.Range("O19").Value = 0.06970509383378017
$ws.Range("S19").Value = 0.1143878462913315
